# Add a new data row (row 3) to the "Artfynd" sheet, mirroring the
# structure of the existing row 2 record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric columns -------------------------------------------------
$ws.Range("A3").Value = 131139952
$ws.Range("B3").Value = 57830
$ws.Range("E3").Value = 100067
$ws.Range("Q3").Value = 411545
$ws.Range("R3").Value = 6578814
$ws.Range("S3").Value = 50

# --- Text columns (plain) --------------------------------------------
$ws.Range("D3").Value = "NT"
$ws.Range("F3").Value = "Havsörn"
$ws.Range("G3").Value = "Haliaeetus albicilla"
$ws.Range("H3").Value = "(Linnaeus, 1758)"

# "Antal" holds the text "1" (not a number) in this sheet, so force the
# cell to text formatting before assigning it.
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "1"

$ws.Range("K3").Value = "adult"
$ws.Range("M3").Value = "förbiflygande"
$ws.Range("P3").Value = "Vidön, Dingelsundet, Vrm"
$ws.Range("T3").Value = "Värmland"
$ws.Range("U3").Value = "Hammarö"
$ws.Range("V3").Value = "Värmland"
$ws.Range("W3").Value = "Hammarö"

# --- Date/time text columns (kept as literal text, not Excel dates) --
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2026-02-12"
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Value = "15:50"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2026-02-12"
$ws.Range("AB3").NumberFormat = "@"
$ws.Range("AB3").Value = "15:50"

# --- Boolean columns ---------------------------------------------------
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false

# --- Reporter / observer columns ---------------------------------------
$ws.Range("AW3").Value = "Peter Adén"
$ws.Range("AX3").Value = "Peter Adén"
